# Added new metadata columns (Background Color, Text Align, Font Bold)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Background Color" column: header + the single value that applies
# to the "Location" field row.
$ws.Range("G1").Value = "Background Color"
$ws.Range("G2").Value = "#FFFFDD"

# New "Text Align" column values for every metadata row.
$ws.Range("H2").Value = "Left"
$ws.Range("I2").Value = "Y"

# Header row for the two right-most new columns.
$ws.Range("H1").Value = "Text Align"
$ws.Range("I1").Value = "Font Bold"

# Remaining rows for "Text Align" / "Font Bold".
$ws.Range("H3").Value = "Left"
$ws.Range("I3").Value = "N"

$ws.Range("H4").Value = "Left"
$ws.Range("I4").Value = "N"

# Size the new column to fit its contents, like Excel does automatically
# when you type into a previously-empty column.
[void]$ws.Columns("G:G").AutoFit()

# Leave the selection where the author ended up after adding the data.
[void]$ws.Range("H14").Select()
